# Weekly price-sheet update: a new daily record is inserted at row 543 of the
# "Hortaliza, Macroferia Regional de Talca - Zapallo italiano" sheet, pushing
# the existing rows 543:613 down to 544:614 (dimension grows from R613 to R614).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 543 (shifts rows 543:613 -> 544:614).
$ws.Rows.Item(543).Insert()

# Populate the newly inserted row 543 with the new record's data.
$ws.Range("A543").Value = 5
$ws.Range("B543").Value = "Macroferia Regional de Talca"
$ws.Range("C543").Value = "Maule"
$ws.Range("D543").Value = 45154
$ws.Range("E543").Value = 7
$ws.Range("F543").Value = 100112032
$ws.Range("G543").Value = "Zapallo italiano"
$ws.Range("H543").Value = "Sin especificar"
$ws.Range("I543").Value = "Primera"
$ws.Range("J543").Value = 300
$ws.Range("K543").Value = 14000
$ws.Range("L543").Value = 14000
$ws.Range("M543").Value = 14000
$ws.Range("N543").Value = "$/caja 50 unidades"
$ws.Range("O543").Value = "Región de Arica y Parinacota"
$ws.Range("P543").Value = 280
$ws.Range("Q543").Value = 50
$ws.Range("R543").Value = "Hortaliza"
